# daily auto push: 2026-01-13 22:36 UTC
# Insert a new day's data row (2026/01/14) into the daily log sheet at row 645,
# pushing every subsequent row (old 645..686) down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 645:686 down to 646:687 by inserting a fresh row at 645.
$ws.Rows(645).Insert()

# Column A stores the date as literal text (e.g. "2026/01/14"), matching the
# rest of the sheet. Force text formatting *before* assigning the value so
# Excel doesn't auto-convert the "yyyy/mm/dd"-shaped string into a real date
# serial number.
$ws.Range("A645").NumberFormat = "@"

$ws.Range("A645").Value = "2026/01/14"
$ws.Range("B645").Value = "水"
$ws.Range("C645").Value = 3
$ws.Range("D645").Value = 201
